# "worked out some bugs" - fill in the rest of row 3 (the next simulated
# trade) and widen column D so the new BuyPrice value isn't clipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (BuyPrice) needs to be a bit wider now that row 3 has a value in it.
$ws.Columns.Item(4).ColumnWidth = 10

# Row 3: Date (A3), BuyPrice (D3) and the Holding flag (G3). C3 (Principle)
# was already populated. Use the same built-in date/time format (numFmtId 22)
# that A2/G2 already carry, so the cells pick up the existing style instead
# of minting a new one.
$ws.Range("A3").Value = 42650.366956018515
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"

$ws.Range("D3").Value = 104.839996

$ws.Range("G3").Value = $true
$ws.Range("G3").NumberFormat = "m/d/yy h:mm"
